# Applies the forecast-shift correction described in the commit:
# "Fixed update to excel issue" - the Week_Start_Date values (and a handful
# of downstream forecast numbers) on the "Forecast Comparison" sheet were
# off by three weeks; the "Summary" sheet's derived figures are refreshed
# to match.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# Column B = Week_Start_Date (stored as plain text in the workbook). A
# leading apostrophe forces Excel to keep the literal text instead of
# re-interpreting the ISO date string as a date serial number - exactly
# what typing a date into a text cell in the Excel UI would do. Columns
# D-H are the numeric forecast figures. Row number -> { column -> new
# value }.

$wsForecast.Range("B2").Value  = "'2025-02-02"

$wsForecast.Range("B3").Value  = "'2025-02-09"

$wsForecast.Range("B4").Value  = "'2025-02-16"

$wsForecast.Range("B5").Value  = "'2025-02-23"

$wsForecast.Range("B6").Value  = "'2025-03-02"
$wsForecast.Range("H6").Value  = 7

$wsForecast.Range("B7").Value  = "'2025-03-09"

$wsForecast.Range("B8").Value  = "'2025-03-16"
$wsForecast.Range("D8").Value  = 2
$wsForecast.Range("E8").Value  = 4
$wsForecast.Range("F8").Value  = 3
$wsForecast.Range("G8").Value  = 6
$wsForecast.Range("H8").Value  = 12

$wsForecast.Range("B9").Value  = "'2025-03-23"
$wsForecast.Range("G9").Value  = 4

$wsForecast.Range("B10").Value = "'2025-03-30"
$wsForecast.Range("E10").Value = 3
$wsForecast.Range("G10").Value = 5
$wsForecast.Range("H10").Value = 10

$wsForecast.Range("B11").Value = "'2025-04-06"
$wsForecast.Range("E11").Value = 3
$wsForecast.Range("G11").Value = 5
$wsForecast.Range("H11").Value = 10

$wsForecast.Range("B12").Value = "'2025-04-13"
$wsForecast.Range("D12").Value = 2
$wsForecast.Range("E12").Value = 4
$wsForecast.Range("F12").Value = 3
$wsForecast.Range("G12").Value = 5
$wsForecast.Range("H12").Value = 11

$wsForecast.Range("B13").Value = "'2025-04-20"
$wsForecast.Range("D13").Value = 2
$wsForecast.Range("H13").Value = 16

$wsForecast.Range("B14").Value = "'2025-04-27"
$wsForecast.Range("D14").Value = 2
$wsForecast.Range("E14").Value = 5
$wsForecast.Range("G14").Value = 7
$wsForecast.Range("H14").Value = 15

$wsForecast.Range("B15").Value = "'2025-05-04"
$wsForecast.Range("H15").Value = 14

$wsForecast.Range("B16").Value = "'2025-05-11"
$wsForecast.Range("H16").Value = 17

$wsForecast.Range("B17").Value = "'2025-05-18"
$wsForecast.Range("D17").Value = 2
$wsForecast.Range("H17").Value = 16

# --- Summary sheet ---------------------------------------------------------
# All "Value" cells on this sheet are plain text, even the ones that hold
# digits only (e.g. "9", "4") - so the numeric-looking "22" needs the same
# force-text treatment as the date strings above.

$wsSummary.Range("B2").Value  = "2023-01-01 to 2025-01-26"
$wsSummary.Range("B9").Value  = "'22"
$wsSummary.Range("B15").Value = "'2025-02-02"
